$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 504899.5
$ws.Range("J48").Value = 504899.5
$ws.Range("L48").Value = 1514698.5
$ws.Range("N48").Value = -1515282.5
$ws.Range("H55").Value = 556.1429
$ws.Range("I55").Value = 482.16666
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 482.16666
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -268.16666
$ws.Range("N55").Value = -1428
$ws.Range("H56").Value = 504899.5
$ws.Range("J56").Value = 504899.5
$ws.Range("L56").Value = 1514698.5
$ws.Range("N56").Value = -1515766.5
$ws.Range("H74").Value = 11978.363
$ws.Range("I74").Value = 11576.2
$ws.Range("J74").Value = 16000
$ws.Range("K74").Value = 11576.2
$ws.Range("L74").Value = 16000
$ws.Range("M74").Value = -10640.2
$ws.Range("N74").Value = -17872
$ws.Range("H77").Value = 11978.363
$ws.Range("I77").Value = 11576.2
$ws.Range("J77").Value = 16000
$ws.Range("K77").Value = 57881
$ws.Range("L77").Value = 80000
$ws.Range("M77").Value = -53201
$ws.Range("N77").Value = -89360
$ws.Range("H96").Value = 1658
$ws.Range("I96").Value = 423
$ws.Range("J96").Value = 4128
$ws.Range("K96").Value = 1269
$ws.Range("L96").Value = 12384
$ws.Range("M96").Value = 104
$ws.Range("N96").Value = -15130
$ws.Range("H100").Value = 3075.1333
$ws.Range("I100").Value = 3223.3572
$ws.Range("K100").Value = 3223.3572
$ws.Range("M100").Value = -2682.3572
$ws.Range("H106").Value = 2657.75
$ws.Range("I106").Value = 2657.75
$ws.Range("K106").Value = 2657.75
$ws.Range("M106").Value = -2026.75
$ws.Range("H135").Value = 793.1667
$ws.Range("I135").Value = 568.26666
$ws.Range("J135").Value = 1917.6666
$ws.Range("K135").Value = 5114.39994
$ws.Range("L135").Value = 17258.9994
$ws.Range("M135").Value = -2579.39994
$ws.Range("N135").Value = -22328.9994
$ws.Range("H138").Value = 6900.875
$ws.Range("J138").Value = 7220.048
$ws.Range("L138").Value = 21660.144
$ws.Range("N138").Value = -31940.144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9575.667
$ws.Range("I32").Value = 7566.95
$ws.Range("K32").Value = 7566.95
$ws.Range("M32").Value = -7279.95
$ws.Range("H45").Value = 2899.5
$ws.Range("I45").Value = 2899
$ws.Range("K45").Value = 2899
$ws.Range("M45").Value = -2522
$ws.Range("H122").Value = 1966
$ws.Range("I122").Value = 1966
$ws.Range("K122").Value = 5898
$ws.Range("M122").Value = -3448
$ws.Range("H132").Value = 4999
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 15426.286
$ws.Range("I99").Value = 17664
$ws.Range("K99").Value = 17664
$ws.Range("M99").Value = -16166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2178.4285
$ws.Range("I31").Value = 1416
$ws.Range("K31").Value = 1416
$ws.Range("M31").Value = -1121
$ws.Range("H34").Value = 2178.4285
$ws.Range("I34").Value = 1416
$ws.Range("K34").Value = 1416
$ws.Range("M34").Value = -1214
$ws.Range("H68").Value = 61855
$ws.Range("J68").Value = 69818.75
$ws.Range("L68").Value = 69818.75
$ws.Range("N68").Value = -71316.75
$ws.Range("H71").Value = 61855
$ws.Range("J71").Value = 69818.75
$ws.Range("L71").Value = 209456.25
$ws.Range("N71").Value = -216944.25
$ws.Range("H74").Value = 70209.336
$ws.Range("J74").Value = 70209.336
$ws.Range("L74").Value = 70209.336
$ws.Range("N74").Value = -71957.336
$ws.Range("H77").Value = 70209.336
$ws.Range("J77").Value = 70209.336
$ws.Range("L77").Value = 210628.008
$ws.Range("N77").Value = -219364.008
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H122").Value = 7996
$ws.Range("I122").Value = 7497
$ws.Range("K122").Value = 22491
$ws.Range("M122").Value = -20041
$ws.Range("H141").Value = 483196.4
$ws.Range("J141").Value = 698660.7
$ws.Range("L141").Value = 698660.7
$ws.Range("N141").Value = -709020.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5249.75
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 30000
$ws.Range("N34").Value = -30168
$ws.Range("H55").Value = 13247.4
$ws.Range("J55").Value = 14163.777
$ws.Range("L55").Value = 42491.331
$ws.Range("N55").Value = -42845.331
$ws.Range("H122").Value = 2065.2856
$ws.Range("I122").Value = 994.5
$ws.Range("J122").Value = 2493.6
$ws.Range("K122").Value = 8950.5
$ws.Range("L122").Value = 22442.4
$ws.Range("M122").Value = -6500.5
$ws.Range("N122").Value = -27342.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2999.25
$ws.Range("I97").Value = 2999.25
$ws.Range("K97").Value = 2999.25
$ws.Range("M97").Value = -2503.25
$ws.Range("H132").Value = 4332.6665
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -20057
$ws.Range("H141").Value = 81999.336
$ws.Range("J141").Value = 81999.336
$ws.Range("L141").Value = 81999.336
$ws.Range("N141").Value = -92359.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H68").Value = 2290
$ws.Range("I68").Value = 2290
$ws.Range("K68").Value = 2290
$ws.Range("M68").Value = -1541
$ws.Range("H71").Value = 2290
$ws.Range("I71").Value = 2290
$ws.Range("K71").Value = 11450
$ws.Range("M71").Value = -7706
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H136").Value = 3853.1538
$ws.Range("I136").Value = 1293
$ws.Range("K136").Value = 3879
$ws.Range("M136").Value = -1329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 60000
$ws.Range("J112").Value = 60000
$ws.Range("L112").Value = 60000
$ws.Range("N112").Value = -62954
$ws.Range("H132").Value = 3100.3809
$ws.Range("I132").Value = 1588.125
$ws.Range("K132").Value = 4764.375
$ws.Range("M132").Value = -2234.375
